# Regenerate save_data: the "K" column (column G) values are recomputed
# (std/mean regen, calc and write s_vals) and rewritten for rows 2-63.
# Everything else in the sheet is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..63 (in row order)
$kValues = @(2,2,1,1,1,0,1,3,1,2,0,3,2,2,1,2,1,2,3,3,0,0,2,2,3,3,1,1,1,3,1,0,1,1,1,0,2,2,2,2,2,2,0,0,2,2,2,0,2,1,0,1,2,5,1,0,2,1,2,3,0,1)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("G$row").Value = $kValues[$i]
}
